$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author replaced curly double quotation marks (U+201C / U+201D) with
# straight single quotation marks (') inside a handful of Kal'tsit English
# (en_US, column C) dialogue lines. Only these four cells changed.

$targets = @("C18", "C23", "C34", "C61")

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    $text = $text.Replace([char]0x201C, "'")
    $text = $text.Replace([char]0x201D, "'")
    $cell.Value2 = $text
}
